{"js": "// Add steps 3-6 of \"MODO DE PREPARO\" after the existing step 2 paragraph,\n// at the end of the document body (immediately before the section\n// properties), matching the unified diff.\nconst body = context.document.body;\n\nconst newSteps = [\n  \"3. Adicione a mistura l\u00edquida \u00e0 tigela e mexa bem. \",\n  \"4. Acrescente o fermento e misture delicadamente. \",\n  \"5. Despeje a massa em uma forma untada e enfarinhada. \",\n  \"6. Asse em forno preaquecido a 180\u00b0C por cerca de 40 minutos.\"\n];\n\nfor (const stepText of newSteps) {\n  body.insertParagraph(stepText, \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Add steps 3-6 of \"MODO DE PREPARO\" after the existing step 2 paragraph,\n# at the end of the document body (immediately before the section\n# properties), matching the unified diff.\n$d = $word.ActiveDocument\n\n$newSteps = @(\n    \"3. Adicione a mistura l\u00edquida \u00e0 tigela e mexa bem. \",\n    \"4. Acrescente o fermento e misture delicadamente. \",\n    \"5. Despeje a massa em uma forma untada e enfarinhada. \",\n    \"6. Asse em forno preaquecido a 180\u00b0C por cerca de 40 minutos.\"\n)\n\nforeach ($stepText in $newSteps) {\n    $lastPara = $d.Paragraphs.Last\n    $lastPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.Text = $stepText\n}\n"}
